# Apply 2022-11-29 daily crime-count update to violent-crime-full-year.xlsx
# Each worksheets column I (year 2022) holds a running year-to-date total;
# some column D (2017) corrections are also included, matching the source diff.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value2 = 6679
$ws.Range("I3").Value2 = 6975
$ws.Range("D4").Value2 = 1940
$ws.Range("I4").Value2 = 1592
$ws.Range("I5").Value2 = 650
$ws.Range("I6").Value2 = 8053
$ws.Range("D7").Value2 = 28130
$ws.Range("I7").Value2 = 23949

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I5").Value2 = 73
$ws.Range("I8").Value2 = 1434
$ws.Range("I9").Value2 = 123
$ws.Range("I10").Value2 = 173
$ws.Range("I12").Value2 = 60
$ws.Range("I14").Value2 = 133
$ws.Range("I18").Value2 = 185
$ws.Range("I19").Value2 = 674
$ws.Range("I20").Value2 = 592
$ws.Range("I24").Value2 = 66
$ws.Range("I25").Value2 = 126
$ws.Range("I26").Value2 = 32
$ws.Range("I27").Value2 = 208
$ws.Range("I29").Value2 = 1444
$ws.Range("I36").Value2 = 328
$ws.Range("D42").Value2 = 1218
$ws.Range("I42").Value2 = 884
$ws.Range("I48").Value2 = 307
$ws.Range("I50").Value2 = 122
$ws.Range("I51").Value2 = 285
$ws.Range("I52").Value2 = 544
$ws.Range("I53").Value2 = 261
$ws.Range("I55").Value2 = 276
$ws.Range("I60").Value2 = 136
$ws.Range("I64").Value2 = 191
$ws.Range("I65").Value2 = 554
$ws.Range("I67").Value2 = 914
$ws.Range("I68").Value2 = 82
$ws.Range("I72").Value2 = 96
$ws.Range("I75").Value2 = 76
$ws.Range("I78").Value2 = 323
$ws.Range("I79").Value2 = 682
$ws.Range("I80").Value2 = 76
$ws.Range("I82").Value2 = 27
$ws.Range("I85").Value2 = 1075
$ws.Range("I87").Value2 = 57
$ws.Range("I91").Value2 = 253
$ws.Range("I92").Value2 = 71
$ws.Range("I93").Value2 = 135
$ws.Range("I94").Value2 = 244
$ws.Range("I95").Value2 = 363
$ws.Range("I96").Value2 = 272
$ws.Range("I97").Value2 = 197
$ws.Range("I99").Value2 = 424
$ws.Range("D101").Value2 = 28130
$ws.Range("I101").Value2 = 23949

# Sheet 3: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value2 = 305
$ws.Range("I3").Value2 = 407
$ws.Range("I4").Value2 = 48
$ws.Range("I6").Value2 = 280
$ws.Range("I7").Value2 = 1075

# Sheet 5: Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I3").Value2 = 177
$ws.Range("I5").Value2 = 19
$ws.Range("I7").Value2 = 544

# Sheet 7: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value2 = 426
$ws.Range("I7").Value2 = 1434

# Sheet 8: Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I6").Value2 = 122
$ws.Range("I7").Value2 = 261

# Sheet 11: West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value2 = 79
$ws.Range("I3").Value2 = 62
$ws.Range("I6").Value2 = 106
$ws.Range("I7").Value2 = 272

# Sheet 12: Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I6").Value2 = 46
$ws.Range("I7").Value2 = 133

# Sheet 15: Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value2 = 119
$ws.Range("I4").Value2 = 28
$ws.Range("I7").Value2 = 424

# Sheet 16: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value2 = 216
$ws.Range("I3").Value2 = 341
$ws.Range("I6").Value2 = 275
$ws.Range("I7").Value2 = 914

# Sheet 19: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("I3").Value2 = 165
$ws.Range("I6").Value2 = 166
$ws.Range("I7").Value2 = 554

# Sheet 21: West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value2 = 126
$ws.Range("I7").Value2 = 363

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value2 = 500
$ws.Range("I5").Value2 = 48
$ws.Range("I6").Value2 = 398
$ws.Range("I7").Value2 = 1444

# Sheet 26: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value2 = 197
$ws.Range("I7").Value2 = 674

# Sheet 28: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value2 = 158
$ws.Range("I7").Value2 = 307

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value2 = 267
$ws.Range("D4").Value2 = 74
$ws.Range("I4").Value2 = 56
$ws.Range("I6").Value2 = 327
$ws.Range("D7").Value2 = 1218
$ws.Range("I7").Value2 = 884

# Sheet 34: Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I2").Value2 = 53
$ws.Range("I3").Value2 = 35
$ws.Range("I7").Value2 = 173

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value2 = 81
$ws.Range("I6").Value2 = 117
$ws.Range("I7").Value2 = 323

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value2 = 85
$ws.Range("I7").Value2 = 276

# Sheet 37: Dunning
$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("I2").Value2 = 24
$ws.Range("I7").Value2 = 66

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value2 = 92
$ws.Range("I7").Value2 = 253

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value2 = 222
$ws.Range("I6").Value2 = 197
$ws.Range("I7").Value2 = 682

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I2").Value2 = 56
$ws.Range("I7").Value2 = 191

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I6").Value2 = 204
$ws.Range("I7").Value2 = 592

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I2").Value2 = 52
$ws.Range("I7").Value2 = 185

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value2 = 94
$ws.Range("I6").Value2 = 103
$ws.Range("I7").Value2 = 328

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I2").Value2 = 38
$ws.Range("I7").Value2 = 135

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I2").Value2 = 46
$ws.Range("I7").Value2 = 244

# Sheet 52: East Side
$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I6").Value2 = 34
$ws.Range("I7").Value2 = 126

# Sheet 56: Lincoln Square
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I2").Value2 = 36
$ws.Range("I3").Value2 = 27
$ws.Range("I7").Value2 = 122

# Sheet 57: East Village
$ws = $wb.Worksheets.Item('East Village')
$ws.Range("I3").Value2 = 6
$ws.Range("I7").Value2 = 32

# Sheet 61: Avalon Park
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I2").Value2 = 42
$ws.Range("I7").Value2 = 123

# Sheet 65: West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I3").Value2 = 34
$ws.Range("I7").Value2 = 197

# Sheet 66: West Elsdon
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I2").Value2 = 22
$ws.Range("I7").Value2 = 71

# Sheet 70: Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I6").Value2 = 29
$ws.Range("I7").Value2 = 73

# Sheet 71: Edgewater
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I6").Value2 = 81
$ws.Range("I7").Value2 = 208

# Sheet 73: Pullman
$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("I6").Value2 = 21
$ws.Range("I7").Value2 = 76

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I6").Value2 = 113
$ws.Range("I7").Value2 = 285

# Sheet 76: North Park
$ws = $wb.Worksheets.Item('North Park')
$ws.Range("I2").Value2 = 27
$ws.Range("I7").Value2 = 82

# Sheet 78: Morgan Park
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I3").Value2 = 35
$ws.Range("I7").Value2 = 136

# Sheet 82: Old Town
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value2 = 47
$ws.Range("I7").Value2 = 96

# Sheet 83: Sheffield & DePaul
$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("I3").Value2 = 8
$ws.Range("I6").Value2 = 27

# Sheet 87: Rush & Division
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I6").Value2 = 42
$ws.Range("I7").Value2 = 76

# Sheet 91: Beverly
$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I6").Value2 = 32
$ws.Range("I7").Value2 = 60

# Sheet 92: Ukrainian Village
$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I3").Value2 = 11
$ws.Range("I7").Value2 = 57
